$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 5 new rows before row 16 (they'll inherit a fresh/blank style) ---
$ws.Range("A16:A20").EntireRow.Insert()

# The insert pushed the old blank rows 16-26 down to 21-31 (old row 28 -> row 33).
# The target layout only keeps 7 of those blank rows (new rows 21-27); the
# trailing 4 blank rows (now at 28-31) are removed entirely.
$ws.Range("A28:A31").EntireRow.Delete()

# --- Fix up formatting of the 5 freshly inserted rows (16-20) so they follow
#     the same alternating style pattern as the rest of the table -------------
$ws.Range("A22:F22").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$ws.Range("A18:F18").PasteSpecial(-4122)
$ws.Range("A20:F20").PasteSpecial(-4122)

$ws.Range("A21:F21").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Range("A19:F19").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# Restore row height for the new rows (matches the rest of the table)
$ws.Range("A16:A20").RowHeight = 15.75

# --- Populate the new component rows -----------------------------------------
# Seed the brand-new shared strings in the same order the source workbook used
# (Regulador 74RM33, SOQUETE torneado 28 pinos, 74HC125 BUS Line Driver) so the
# shared-string table indices line up.
$ws.Range("B19").Value = "Regulador 74RM33"
$ws.Range("B18").Value = "SOQUETE torneado 28 pinos"
$ws.Range("B16").Value = "74HC125 BUS Line Driver"

# Row 16: Regulador datasheet component (74HC125 BUS Line Driver)
$ws.Range("A16").Value = 41239
$ws.Range("C16").Value = 1.5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = "Mundial Componentes"
$ws.Range("F16").Formula = "=PRODUCT(C16:D16)"

# Row 17
$ws.Range("A17").Value = 41243
$ws.Range("B17").Value = "CONECTOR modu 2542"
$ws.Range("C17").Value = 0.6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = "Mundial Componentes"
$ws.Range("F17").Formula = "=PRODUCT(C17:D17)"

# Row 18
$ws.Range("A18").Value = 41243
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "Mundial Componentes"
$ws.Range("F18").Formula = "=PRODUCT(C18:D18)"

# Row 19
$ws.Range("A19").Value = 41243
$ws.Range("C19").Value = 2.5
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = "Mundial Componentes"
$ws.Range("F19").Formula = "=PRODUCT(C19:D19)"

# Row 20
$ws.Range("A20").Value = 41243
$ws.Range("B20").Value = "TERMINAL para conector modu 22/26AWG"
$ws.Range("C20").Value = 0.1
$ws.Range("D20").Value = 16
$ws.Range("E20").Value = "Mundial Componentes"
$ws.Range("F20").Formula = "=PRODUCT(C20:D20)"

# --- View state tweaks from the diff -----------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("C17").Select()
